# "Subida luego del webinario"
#
# 1. Bump the cached datetimeFigureOut date field (29/10/2020 -> 30/10/2020)
#    on the Slide Master and every Slide Layout (12 occurrences total).
# 2. On slide 2 ("Estructura del webinario"), shrink the bullet list back
#    down (it no longer needs the extra autofit shrink), and drop the
#    "Grupo de R" / "Grupo de Python" bullets (and their matching
#    click-to-reveal animations) that are no longer part of the agenda.

$p = $ppt.ActivePresentation

# --- 1. Date placeholder on the Slide Master ---------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "29/10/2020") {
        $sh.TextFrame.TextRange.Text = "30/10/2020"
    }
}

# --- 1b. Date placeholder on every Slide Layout -------------------------
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "29/10/2020") {
            $sh.TextFrame.TextRange.Text = "30/10/2020"
        }
    }
}

# --- 2. Slide 2 content placeholder -------------------------------------
$slide2 = $p.Slides.Item(2)
$content = $slide2.Shapes.Item(2)

# Remove the two trailing bullets ("Grupo de R" / "Grupo de Python").
# They sit at paragraphs 9 and 10 (1-based); after the first delete the
# next one slides up into slot 9, so deleting slot 9 twice removes both.
$tr = $content.TextFrame.TextRange
$tr.Paragraphs(9, 1).Delete()
$tr.Paragraphs(9, 1).Delete()

# The placeholder no longer needs the manual line-spacing shrink, and it
# moves up / shrinks a bit now that it has two fewer lines.
$content.Top = 130.789
$content.Height = 358.0184
$content.TextFrame.AutoSize = 2

# Drop the two click-to-reveal animations that targeted the removed
# bullets (they were the last two effects in the main sequence).
$seq = $slide2.TimeLine.MainSequence
$seq.Item($seq.Count).Delete()
$seq.Item($seq.Count).Delete()
